$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new quarter labels for S1:U1, matching the
# formatting (bold, centered, bordered) already used by the other header cells ---
$ws.Range("R1").Copy()
$ws.Range("S1:U1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("S1").Value = "31/12/2023"
$ws.Range("T1").Value = "31/03/2024"
$ws.Range("U1").Value = "30/06/2024"

# --- Blank separator rows (section headers with no numeric data) ---
# These rows only have label text in column A; columns B:R are empty
# placeholder cells. Copy that same "empty" formatting into S:U so the
# row stays structurally consistent.
$blankRows = @(57, 58, 71, 72, 73, 77, 78)
foreach ($r in $blankRows) {
    $ws.Range("B$r").Copy()
    $ws.Range("S${r}:U${r}").PasteSpecial(-4122)  # xlPasteFormats
}

# --- Data rows: three new quarterly columns (S = 31/12/2023,
# T = 31/03/2024, U = 30/06/2024) appended after the existing R column ---
$data = @"
2,42747015.168,41646391.296,40759304.192
3,15987221.504,15182557.184,13934275.584
4,3750944,4190643.968,1596338.048
5,4024056.064,1665371.008,1918401.024
6,3547310.08,3698857.984,4524867.072
7,3087395.072,3709506.048,3982482.944
8,0,0,0
9,784092.992,952288,921057.9840000001
10,0,0,0
11,793424,965889.9840000001,991129.024
12,5681544.192,5311828.992,4623170.048
13,0,0,0
14,0,0,0
15,0,0,0
16,806582.0159999999,332494.016,375691.008
17,0,0,0
18,0,0,0
19,2200695.04,2421377.024,1681687.04
20,0,0,0
21,0,0,0
22,0,0,0
23,4508383.232,4421430.784,4629020.16
24,16569866.24,16730576.896,17572839.424
25,0,0,0
26,42747015.168,41646391.296,40759304.192
27,10413499.392,10296062.976,9476476.927999999
28,1019688,1138552.96,954649.024
29,5302477.824,5350676.992,5757737.984
30,1543202.048,818616,746515.968
31,163844,128090,92038
32,0,0,0
33,294231.008,1018288,41442
34,1598754.976,1248803.968,1374825.008
35,491300.992,593035.008,509268.992
36,0,0,0
37,9213050.880000001,9262936.063999999,9027227.648
38,5947857.92,6059054.08,5653753.856
39,0,0,0
40,1375278.976,1353860.992,1424324.992
41,328089.984,361124,404880
42,0,0,0
43,1561824,1488897.024,1544269.056
44,0,0,0
45,0,0,0
46,17226,16715,17925
47,23103237.872,22070680.328,22237674.616
48,12484514.816,12484514.816,12484514.816
49,10302265.344,10336512,10405700.608
50,0,0,0
51,780307.968,51703,50265
52,0,-934816,-1793730.048
53,-463848.992,132766,1090924.032
54,0,0,0
55,0,0,0
56,0,0,0
59,4126319.104,6105252.864,7352631.808
60,-1881811.2,-2127079.936,-2610971.904
61,2244507.648,3978172.928,4741660.16
62,-1264030.976,-2497074.944,-3012726.016
63,-167604.992,-1043043.008,-1207831.04
64,-145266.992,-235884,-137284
65,0,0,0
66,-669122.944,-46936,-119993
67,0,0,0
68,-570177.92,-361215.008,-135324.992
69,-2860289.024,0,0
70,2290111.232,-361215.008,-135324.992
74,-571695.872,-205980,128501
75,0,0,0
76,0,0,0
79,-1000,310,176
80,-2666419.968,-934816,-858913.9840000001
"@

$lines = $data -split "`n"
foreach ($line in $lines) {
    $line = $line.Trim()
    if ($line -eq "") { continue }
    $parts = $line -split ","
    $r = $parts[0]
    $sVal = [double]$parts[1]
    $tVal = [double]$parts[2]
    $uVal = [double]$parts[3]
    $ws.Range("S$r").Value = $sVal
    $ws.Range("T$r").Value = $tVal
    $ws.Range("U$r").Value = $uVal
}

Write-Host "Added columns S:U (31/12/2023, 31/03/2024, 30/06/2024) across all 80 rows."
